# Atualizacao de bases das ligas - Switzerland Challenge League
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap full data rows (columns B:AC), keeping column A (row index/id) untouched ---
function Swap-Rows([int]$r1, [int]$r2) {
    $rangeA = $ws.Range("B$r1`:AC$r1")
    $rangeB = $ws.Range("B$r2`:AC$r2")
    $valsA = $rangeA.Value()
    $valsB = $rangeB.Value()
    $rangeA.Value = $valsB
    $rangeB.Value = $valsA
}

Swap-Rows 39 40
Swap-Rows 72 73
Swap-Rows 90 91
Swap-Rows 105 106
Swap-Rows 108 109

# --- Row 132: update with new match data (FC Vaduz vs FC Thun, played 2024-04-01) ---
$ws.Range("B132").Value = 7617793
$ws.Range("E132").Value = 45383.38541666666
$ws.Range("F132").Value = "FC Vaduz"
$ws.Range("G132").Value = "FC Thun"
$ws.Range("K132").Value = 3.6
$ws.Range("L132").Value = 3.75
$ws.Range("M132").Value = 1.909
$ws.Range("N132").Value = 3.8
$ws.Range("O132").Value = 4
$ws.Range("P132").Value = 1.85
$ws.Range("Q132").Value = 0.5
$ws.Range("R132").Value = 1.95
$ws.Range("S132").Value = 1.85
$ws.Range("T132").Value = 3
$ws.Range("U132").Value = 1.875
$ws.Range("V132").Value = 1.925

# --- Remove the now-duplicate row 133 (its data was merged into row 132 above) ---
$ws.Rows.Item(133).Delete()
